$p = $ppt.ActivePresentation
$s = $p.Slides.Add(1, [Microsoft.Office.Interop.PowerPoint.PpSlideLayout]::ppLayoutTitle)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "sdassdaads"
